$wb = $excel.ActiveWorkbook

# --- Summary sheet updates ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.99
$summary.Range("B4").Value = -0.01
$summary.Range("B5").Value = -0.04
$summary.Range("B6").Value = 5
$summary.Range("B8").Value = 2
$summary.Range("B9").Value = 40

# --- Strategy Status sheet updates (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.98999999999999
$status.Range("D4").Value = 5
$status.Range("E4").Value = -0.01
$status.Range("F4").Value = -0.01
$status.Range("G4").Value = 40

# --- New trade row (#5) appended to both "All Trades" and "MarketMaking" sheets ---
$tradeRow = @(5, "2026-02-17", "04:06:14", "MarketMaking", "DOWN", 0.8100000000000001, 0.78, "CLOSED", -3.7037, -0.03, 99.98999999999999, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.11)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 6
    for ($col = 1; $col -le $tradeRow.Length; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        if ($col -eq 2) {
            # Date column: force text so "2026-02-17" isn't auto-converted
            # to a date serial number, then strip the temporary format so
            # no stray style is left behind on the cell.
            $cell.NumberFormat = "@"
            $cell.Value = $tradeRow[$col - 1]
            $cell.ClearFormats()
        } else {
            $cell.Value = $tradeRow[$col - 1]
        }
    }
}
